$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row columns: "<name>_old" -> "<name>_FV2410" and
# "<name>_new" -> "<name>_FV2504" (columns A-J are the "old"/FV2410 side,
# L-U are the "new"/FV2504 side; column K is the unchanged "diff" column).
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    $cell.Value = $val.Replace("_old", "_FV2410")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    $cell.Value = $val.Replace("_new", "_FV2504")
}

# Turn the data range into an Excel Table (ListObject) with autofilter,
# using the renamed header row as the column names.
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
